# Add a new row to the "Completed" reading-list sheet for the book
# "When The Air Hits Your Brain" by Frank Vertosick.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 85

# Copy date formatting (style) from the row above so the new date cells
# reuse the existing date number format instead of minting a new one.
$ws.Range("C84").Copy()
$ws.Range("C$row").PasteSpecial(-4122)
$ws.Range("D84").Copy()
$ws.Range("D$row").PasteSpecial(-4122)

$ws.Range("A$row").Value = "When The Air Hits Your Brain"
$ws.Range("B$row").Value = "Frank Vertosick"
$ws.Range("C$row").Value = Get-Date -Year 2020 -Month 6 -Day 6 -Hour 0 -Minute 0 -Second 0
$ws.Range("D$row").Value = Get-Date -Year 2020 -Month 6 -Day 8 -Hour 0 -Minute 0 -Second 0
$ws.Range("E$row").Value = "medicine;neurology;neuro surgery;surgery;doctor;residency"
$ws.Range("F$row").Value = "Audio"
$ws.Range("G$row").Value = "8 Hours 43 Mins"

$ws.Range("A86").Select()
